$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Wnt2/Fzd7 -> ECs)
$ws.Range("G2").Value = 0.6811249999999999
$ws.Range("M2").Value = 1.123319
$ws.Range("N2").Value = 3.369957
$ws.Range("O2").Value = 0.05053686506648315
$ws.Range("P2").Value = 0.05053686506648315
$ws.Range("Q2").Value = 0.7651206538749999
$ws.Range("R2").Value = 6.886085884875
$ws.Range("S2").Value = 0.05053686506648315
$ws.Range("T2").Value = 0.05053686506648315

# Row 3 (FAPs -> Wnt2/Fzd7 -> FAPs)
$ws.Range("G3").Value = 0.6811249999999999
$ws.Range("O3").Value = 0.5042195746532222
$ws.Range("P3").Value = 0.5042195746532223
$ws.Range("Q3").Value = 7.633809698083331
$ws.Range("R3").Value = 68.70428728274999
$ws.Range("S3").Value = 0.5042195746532222
$ws.Range("T3").Value = 0.5042195746532223

# Row 4 (FAPs -> Wnt2/Fzd7 -> MuSCs)
$ws.Range("G4").Value = 0.6811249999999999
$ws.Range("M4").Value = 4.958620666666667
$ws.Range("N4").Value = 14.875862
$ws.Range("O4").Value = 0.2230827962023326
$ws.Range("P4").Value = 0.2230827962023326
$ws.Range("Q4").Value = 3.377440501583333
$ws.Range("R4").Value = 30.39696451425
$ws.Range("S4").Value = 0.2230827962023326
$ws.Range("T4").Value = 0.2230827962023326

# Row 5 (FAPs -> Wnt2/Fzd7 -> Resolving-Mac)
$ws.Range("G5").Value = 0.6811249999999999
$ws.Range("M5").Value = 4.938126
$ws.Range("N5").Value = 14.814378
$ws.Range("O5").Value = 0.222160764077962
$ws.Range("P5").Value = 0.222160764077962
$ws.Range("Q5").Value = 3.36348107175
$ws.Range("R5").Value = 30.27132964575
$ws.Range("S5").Value = 0.222160764077962
$ws.Range("T5").Value = 0.222160764077962
